$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.801.64'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '3.100.52'
$ws.Range('E3').Value = '  +3.88%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '387.92'
$ws.Range('E5').Value = '  +1.68%  '
$ws.Range('D6').Value = '103.72'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').Value = '0.545'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '0.590'
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('D10').Value = '37.11'
$ws.Range('E10').Value = '  +1.22%  '
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').Value = '0.0865'
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('D13').Value = '3.594.37'
$ws.Range('E13').Value = '  +3.99%  '
$ws.Range('D14').Value = '18.75'
$ws.Range('E14').Value = '  +1.46%  '
$ws.Range('D15').Value = '7.85'
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('D16').Value = '3.099.70'
$ws.Range('E16').Value = '  +4.08%  '
$ws.Range('D17').Value = '0.983'
$ws.Range('E17').Value = '  -1.28%  '
$ws.Range('D18').Value = '10.72'
$ws.Range('E18').Value = '  -4.66%  '
$ws.Range('D19').Value = '51.946.86'
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('E20').Value = '  +2.58%  '
$ws.Range('D21').Value = '12.55'
$ws.Range('E21').Value = '  -0.37%  '
$ws.Range('D22').Value = '0.0₃0971'
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('D23').Value = '70.33'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = '269.25'
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('D25').Value = '3.15'
$ws.Range('E25').Value = '  -2.82%  '
$ws.Range('D26').Value = '8.21'
$ws.Range('E26').Value = '  +4.94%  '
$ws.Range('D27').Value = '27.09'
$ws.Range('E27').Value = '  +3.82%  '
$ws.Range('E28').Value = '  +2.44%  '
$ws.Range('D29').Value = '7.26'
$ws.Range('E29').Value = '  -1.75%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  -0.49%  '
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('D33').Value = '35.66'
$ws.Range('E33').Value = '  +2.46%  '
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('D35').Value = '50.36'
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('E36').Value = '  +1.92%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').Value = '3.39'
$ws.Range('E38').Value = '  +3.19%  '
$ws.Range('D39').Value = '0.297'
$ws.Range('E39').Value = '  +9.15%  '
$ws.Range('E40').Value = '  +2.45%  '
$ws.Range('D41').Value = '17.05'
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').Value = '2.59'
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').Value = '127.57'
$ws.Range('E43').Value = '  +3.38%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '0.116'
$ws.Range('E44').Value = '  -0.73%  '
$ws.Range('D45').Value = '3.71'
$ws.Range('E45').Value = '  -2.26%  '
$ws.Range('D46').Value = '22.13'
$ws.Range('E46').Value = '  +3.28%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '2.47'
$ws.Range('E47').Value = '  +4.69%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '2.09'
$ws.Range('E48').Value = '  +2.66%  '
$ws.Range('D49').Value = '2.052.61'
$ws.Range('E49').Value = '  +1.19%  '
$ws.Range('D50').Value = '3.414.14'
$ws.Range('E50').Value = '  +4.01%  '
$ws.Range('E51').Value = '  +7.11%  '
